$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held the plain header numbers 1..8 above the "Lambda" row; the new
# layout no longer shows that header row, so just wipe its contents (the
# row itself is not deleted - rows 3-8 below keep their original numbers).
$ws.Rows("2:2").ClearContents()

# The former "Lambda" row (C3) becomes the first row of the SRPT results
# block; rename its label accordingly.
$ws.Range("C3").Value2 = "SRPT"

# Add a brand-new results block below for the RANDOM (no-preemption)
# policy. Row 11 reuses the same Lambda values as row 3 (copied so the
# underlying numeric values match exactly), and rows 12-16 carry the same
# metric labels as rows 4-8, ready to be filled in later.
$ws.Range("D3:K3").Copy($ws.Range("D11:K11"))
$ws.Range("C11").Value2 = "RANDOM (senza prelazione)"

$ws.Range("C12").Value2 = "utilizzo"
$ws.Range("C13").Value2 = "thro"
$ws.Range("C14").Value2 = "med resp"
$ws.Range("C15").Value2 = "med wait"
$ws.Range("C16").Value2 = "max q"

$null = $ws.Range("D12").Select()

# Reposition charts to match the new layout (charts moved down/right to make
# room for the added data block).
$co2 = $ws.ChartObjects("Grafico 2")
$co2.Left = 804.382062
$co2.Top = 465.444488
$co2.Width = 499.944488
$co2.Height = 417.555591

$co3 = $ws.ChartObjects("Grafico 3")
$co3.Left = 0.000000
$co3.Top = 358.333228
$co3.Width = 592.854266
$co3.Height = 467.333386

$co4 = $ws.ChartObjects("Grafico 4")
$co4.Left = 660.180585
$co4.Top = 940.999921
$co4.Width = 646.145965
$co4.Height = 466.222283

$co5 = $ws.ChartObjects("Grafico 5")
$co5.Left = 2.222283
$co5.Top = 902.333386
$co5.Width = 545.527869
$co5.Height = 533.555354
